$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor values (ax, ay, az, gx, gy, gz) for rows 2-21 (columns C-H)
$data = @{
    2  = @(1.136592060327531, -0.5052947402000403, 1.006910741329194, -5.196477890014648, -3.036305665969849, 4.782464027404785)
    3  = @(-3.102525949478149, -7.741808414459229, 0.0443248748779296, -5.916993618011475, -3.426343202590942, -1.53388261795044)
    4  = @(-0.8462124168872824, -1.410848729312418, -6.25395178794861, 0.4529565870761871, 0.4115704596042633, -1.645212769508362)
    5  = @(4.942803740501414, -7.510328069329272, -1.564183712005606, 3.627165079116821, 3.22368860244751, -1.423163294792175)
    6  = @(3.684623420238498, -6.547261834144596, 2.809367418289182, 1.755779385566711, 0.2069306075572967, -0.0197004042565822)
    7  = @(-2.489120721817007, -3.182153344154361, 2.89736366271973, 1.32527220249176, 3.478113889694214, 4.531856536865234)
    8  = @(-2.860559403896316, -2.504158109426502, 1.881697505712511, -0.9728028774261476, 3.588375091552734, -0.8292495608329773)
    9  = @(2.841274738311764, -2.992074728012082, 1.633560657501217, 0.9275988936424256, 3.912744522094727, -0.8049677014350891)
    10 = @(0.4830425977706894, -0.9799425601959225, -0.1168481409549758, 0.65958172082901, -0.1707368344068527, 0.1705841124057769)
    11 = @(2.629005432128906, -1.44602632522583, 6.930900573730469, -0.0062613687478005, -1.035263895988464, 2.51661205291748)
    12 = @(0.08325040340423584, -1.20862039923668, 1.943817764520645, -0.0650571510195732, -0.7725918292999268, -0.6349944472312927)
    13 = @(-1.441726684570314, -0.9328206181526182, 5.954500854015357, 1.70736837387085, -1.938886284828186, -2.866790533065796)
    14 = @(-0.886786460876462, 0.562802791595463, 2.521360546350466, 0.0068722339347004, -3.644727468490601, -0.365755558013916)
    15 = @(0.6777331829071001, -1.399440765380847, 0.4171198606491046, -0.09926560521125791, -5.367214679718018, -0.7640396952629089)
    16 = @(1.204387292265894, -3.619132399559024, -0.313209235668175, 0.1264491081237793, -3.759417533874512, -0.6985244154930115)
    17 = @(1.067639499902725, -3.350790739059448, -3.138641357421875, 0.3510948121547699, -1.851074457168579, -0.6296493411064148)
    18 = @(3.527233093976973, -5.226221919059754, -3.07070302963257, -0.2174680233001709, 1.732872009277344, -2.470797300338745)
    19 = @(7.073627471923828, -3.28519868850708, -0.3319654464721679, 1.422094345092773, 4.435187339782715, 0.8324565887451172)
    20 = @(2.874565660953519, -0.3619521260261513, -1.012392401695252, 2.657721996307373, 5.389358520507812, 2.652835130691528)
    21 = @(-6.492176651954649, -0.7134745419025422, 7.389736890792848, 0.5099197626113892, 0.4361577928066253, 2.850908041000366)
}

$columns = @("C", "D", "E", "F", "G", "H")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}

# Remove the last row (row 22), which shifts the data window up and
# shrinks the used range to A1:H21
$ws.Rows.Item(22).Delete()
